$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column B holds a date-like string ("2024-06-12"). Excel would normally
# auto-convert that to a real date serial on entry, but the source data
# keeps it as plain text (same as every other row), so pre-format the
# cell as text before assigning, then clear the formatting again
# afterwards so the new row ends up with no explicit style -- matching
# the rest of the sheet's unstyled data rows.
$ws.Cells.Item(23, 2).NumberFormat = "@"

$ws.Cells.Item(23, 1).Value = "segqua0708"
$ws.Cells.Item(23, 2).Value = "2024-06-12"
$ws.Cells.Item(23, 3).Value = "julio"
$ws.Cells.Item(23, 4).Value = "Manel, Bernardo"

$ws.Range("A23:D23").ClearFormats()
